# Fruta / hortaliza, semanal
#
# A new weekly price record (dated 2022-07-27) for "Agrícola del Norte S.A.
# de Arica" / Maracuyá / Primera is inserted above the existing row 132,
# pushing the previous rows 132-136 down to 133-137.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 132, shifting the existing rows 132:136 down to 133:137.
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(132, 1).Value = 1
$ws.Cells.Item(132, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(132, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(132, 4).Value = 44769
$ws.Cells.Item(132, 5).Value = 15
$ws.Cells.Item(132, 6).Value = 'Fruta'
$ws.Cells.Item(132, 7).Value = 100108
$ws.Cells.Item(132, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(132, 9).Value = 100108003
$ws.Cells.Item(132, 10).Value = 'Maracuyá'
$ws.Cells.Item(132, 11).Value = 'Sin especificar'
$ws.Cells.Item(132, 12).Value = 'Primera'
$ws.Cells.Item(132, 13).Value = 100
$ws.Cells.Item(132, 14).Value = 24000
$ws.Cells.Item(132, 15).Value = 25000
$ws.Cells.Item(132, 16).Value = 24500
$ws.Cells.Item(132, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(132, 18).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(132, 19).Value = 1225
$ws.Cells.Item(132, 20).Value = 20
